$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Split the "Purpose" paragraph's sentence: keep the opening clause,
#    replace the remainder with the new "work done to research..." text.
# ---------------------------------------------------------------------
$quote = [char]0x201C
$endquote = [char]0x201D

$findRange = $d.Content
$findRange.Find.ClearFormatting()
$findRange.Find.Execute("outline the ") | Out-Null

$p3 = $d.Paragraphs.Item(3)
$tailRange = $d.Range($findRange.End, $p3.Range.End - 1)

$newTail = 'work done to research and find a solution to how the ' + $quote + 'Attendance Made Easy' + $endquote + ' website will be hosted. This is so that appropriate story points can be awarded and distributed to the appropriate person.'

$tailXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>' + $newTail + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$tailRange.InsertXML($tailXml)

# ---------------------------------------------------------------------
# 2. Insert the new "Findings" section (heading + two narrative
#    paragraphs) before the final (bookmarked) paragraph. An extra
#    leading "<w:p/>" becomes the blank divider paragraph right after
#    the Purpose text; the trailing narrative paragraph is left
#    "open" so its runs land inside the existing bookmark paragraph,
#    right before the "_GoBack" bookmark.
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$insertPoint = $d.Range($p4.Range.Start, $p4.Range.Start)

$body = '<w:p/>'
$body += '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Findings</w:t></w:r></w:p>'
$body += '<w:p><w:r><w:t xml:space="preserve">A team member reached out to Dr. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Droz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, a professor at Louisiana Tech who specializes in cloud computing and talked about possible ways to host the website. He recommended using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TechXplore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, a server on the Louisiana Tech campus, to host the website.</w:t></w:r></w:p>'
$body += '<w:p><w:r><w:t xml:space="preserve">He outlined setting up a Virtual Machine (VM) on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TechXplore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and using that to setup LAMP/WAMP on said VM. This will allow the website to be hosted and accessed via an IP address. Then, Dr. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Droz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> mentioned using an online DNS rerouting service, such as </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>DuckDNS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, to reroute traffic so the website can have a custom domain.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($xml)
